# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values, updated per latest spot price data refresh
$ws.Range("A2").Value = 46068
$ws.Range("B2").Value = 4.83
$ws.Range("C2").Value = 1.64
$ws.Range("D2").Value = 0.24
$ws.Range("E2").Value = 0.01
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.01
$ws.Range("I2").Value = 0.06
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = -0.63
$ws.Range("L2").Value = -1.19
$ws.Range("M2").Value = -1.6
$ws.Range("N2").Value = -1.6
$ws.Range("O2").Value = -1.83
$ws.Range("P2").Value = -2.2
$ws.Range("Q2").Value = -2.41
$ws.Range("R2").Value = -1.7
$ws.Range("S2").Value = -0.47
$ws.Range("T2").Value = 0.13
$ws.Range("U2").Value = 3.4
$ws.Range("V2").Value = 9.449999999999999
$ws.Range("W2").Value = 9.25
$ws.Range("X2").Value = 6.07
$ws.Range("Y2").Value = 0.35
$ws.Range("Z2").Value = 0.91
$ws.Range("AB2").Value = 6.28
$ws.Range("AD2").Value = 9.35
$ws.Range("AE2").Value = "0h-2h"
$ws.Range("AF2").Value = 3.24
$ws.Range("AG2").Value = "2h-23h"
